$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell reference -> new text value, derived from the commit diff.
# Values are percentages / prices stored as literal text (inlineStr) in the
# source sheet, so each cell is pre-formatted as Text ("@") before the
# assignment to stop Excel from auto-coercing the numeric-looking strings
# (e.g. "1.61%" or "305.62") into Number/Percentage cells.
$updates = [ordered]@{
    "D2" = "305.62"
    "E2" = "1.61%"
    "D3" = "36.34"
    "E3" = "-3.32%"
    "E4" = "1.20%"
    "E5" = "0.08%"
    "D6" = "2.174"
    "E6" = "-1.12%"
    "D7" = "7.920"
    "E7" = "-1.07%"
    "D8" = "0.9188"
    "E8" = "0.78%"
    "D9" = "0.09738"
    "E9" = "5.60%"
    "E10" = "-0.26%"
    "D11" = "0.08640"
    "E11" = "1.92%"
    "D12" = "0.03486"
    "E12" = "-1.45%"
    "D13" = "0.09930"
    "E13" = "-0.07%"
    "D14" = "0.001446"
    "E14" = "-1.59%"
    "D15" = "0.005630"
    "E15" = "-0.25%"
    "D16" = "3.461"
    "E16" = "-0.46%"
    "D17" = "4.094"
    "E17" = "2.52%"
    "D18" = "2.392"
    "E18" = "14.16%"
    "D20" = "0.1361"
    "E20" = "3.67%"
    "D21" = "4.780"
    "E21" = "5.02%"
    "E22" = "-1.65%"
    "D23" = "0.04563"
    "E23" = "-1.50%"
    "D24" = "0.005097"
    "E24" = "14.73%"
    "E25" = "0.37%"
    "D26" = "0.0001401"
    "E26" = "7.82%"
    "D27" = "0.0004753"
    "E27" = "0.14%"
    "D39" = "0.01829"
    "E39" = "4.66%"
    "D40" = "0.04749"
    "E40" = "0.96%"
    "D41" = "0.007780"
    "E41" = "-1.59%"
    "D42" = "0.1398"
    "E42" = "0.48%"
    "D43" = "0.007732"
    "E43" = "0.97%"
    "D44" = "0.002239"
    "E44" = "-2.61%"
    "D45" = "0.01133"
    "E45" = "10.71%"
    "D46" = "0.00006319"
    "E46" = "4.32%"
    "D47" = "0.00000000751"
    "E47" = "0.15%"
    "D48" = "0.0005802"
    "E48" = "0.02%"
    "E49" = "176.57%"
    "D50" = "0.002001"
    "E50" = "-25.84%"
    "D51" = "0.00002102"
    "E51" = "0.15%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}

Write-Output "Updated $($updates.Count) cells"
